$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (REG VAL)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.07378995"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.16836306"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 72
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 2

# Row 5 (REG TEST)
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0.15078978"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "0.2989677"
$ws.Range("C5").Style = "Normal"
